$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item(1)
$wsARM = $wb.Worksheets.Item(2)
$wsBSM = $wb.Worksheets.Item(3)
$wsCRP = $wb.Worksheets.Item(4)
$wsCUL = $wb.Worksheets.Item(5)
$wsGSM = $wb.Worksheets.Item(6)
$wsLTW = $wb.Worksheets.Item(7)
$wsWVR = $wb.Worksheets.Item(8)

$wsALC.Range("H2").Value = 94.5
$wsALC.Range("I2").Value = 83.40000000000001
$wsALC.Range("K2").Value = 83.40000000000001
$wsALC.Range("M2").Value = 29.59999999999999

$wsALC.Range("H4").Value = 419.2
$wsALC.Range("I4").Value = 419.2
$wsALC.Range("K4").Value = 419.2
$wsALC.Range("M4").Value = -305.2

$wsALC.Range("H6").Value = 2854.4285
$wsALC.Range("I6").Value = 2854.4285
$wsALC.Range("K6").Value = 8563.2855
$wsALC.Range("M6").Value = -8451.2855

$wsALC.Range("H9").Value = 227.44444
$wsALC.Range("I9").Value = 130.6
$wsALC.Range("J9").Value = 348.5
$wsALC.Range("K9").Value = 130.6
$wsALC.Range("L9").Value = 348.5
$wsALC.Range("M9").Value = 38.40000000000001
$wsALC.Range("N9").Value = -686.5

$wsALC.Range("H48").Value = 10000
$wsALC.Range("J48").Value = 10000
$wsALC.Range("L48").Value = 30000
$wsALC.Range("N48").Value = -30584

$wsALC.Range("H56").Value = 10000
$wsALC.Range("J56").Value = 10000
$wsALC.Range("L56").Value = 30000
$wsALC.Range("N56").Value = -31068

$wsALC.Range("H86").Value = 153848960
$wsALC.Range("J86").Value = 200003970
$wsALC.Range("L86").Value = 200003970
$wsALC.Range("N86").Value = -200006216

$wsALC.Range("H88").Value = 10102526
$wsALC.Range("J88").Value = 1525.5714
$wsALC.Range("L88").Value = 1525.5714
$wsALC.Range("N88").Value = -2337.5714

$wsALC.Range("H89").Value = 153848960
$wsALC.Range("J89").Value = 200003970
$wsALC.Range("L89").Value = 1000019850
$wsALC.Range("N89").Value = -1000031082

$wsALC.Range("H91").Value = 10102526
$wsALC.Range("J91").Value = 1525.5714
$wsALC.Range("L91").Value = 1525.5714
$wsALC.Range("N91").Value = -4333.5714

$wsALC.Range("H98").Value = 1751.9584
$wsALC.Range("I98").Value = 1523.7826
$wsALC.Range("K98").Value = 1523.7826
$wsALC.Range("M98").Value = -25.7826

$wsALC.Range("H111").Value = 1370.909
$wsALC.Range("I111").Value = 1261.25
$wsALC.Range("J111").Value = 1663.3334
$wsALC.Range("K111").Value = 3783.75
$wsALC.Range("L111").Value = 4990.0002
$wsALC.Range("M111").Value = -716.75
$wsALC.Range("N111").Value = -11124.0002

$wsALC.Range("H122").Value = 1751.9584
$wsALC.Range("I122").Value = 1523.7826
$wsALC.Range("K122").Value = 4571.3478
$wsALC.Range("M122").Value = -2121.3478

$wsALC.Range("H137").Value = 139048.33
$wsALC.Range("I137").Value = 207050
$wsALC.Range("K137").Value = 621150
$wsALC.Range("M137").Value = -618600

$wsARM.Range("H32").Value = 3693.2222
$wsARM.Range("I32").Value = 1793.125
$wsARM.Range("K32").Value = 1793.125
$wsARM.Range("M32").Value = -1506.125

$wsARM.Range("H45").Value = 3834.8462
$wsARM.Range("I45").Value = 3255.35
$wsARM.Range("K45").Value = 3255.35
$wsARM.Range("M45").Value = -2878.35

$wsARM.Range("H46").Value = 0
$wsARM.Range("J46").Value = 0
$wsARM.Range("L46").Value = 0
$wsARM.Range("N46").ClearContents()

$wsARM.Range("H102").Value = 1662.4166
$wsARM.Range("I102").Value = 1244.9
$wsARM.Range("K102").Value = 1244.9
$wsARM.Range("M102").Value = 377.0999999999999

$wsARM.Range("H122").Value = 4146.421
$wsARM.Range("I122").Value = 2610.88
$wsARM.Range("K122").Value = 7832.64
$wsARM.Range("M122").Value = -5382.64

$wsARM.Range("H132").Value = 2599.34
$wsARM.Range("I132").Value = 2327.125
$wsARM.Range("K132").Value = 6981.375
$wsARM.Range("M132").Value = -4451.375

$wsARM.Range("H133").Value = 93805
$wsARM.Range("J133").Value = 93805
$wsARM.Range("L133").Value = 93805
$wsARM.Range("N133").Value = -98865

$wsARM.Range("H135").Value = 44705.8
$wsARM.Range("J135").Value = 44705.8
$wsARM.Range("L135").Value = 44705.8
$wsARM.Range("N135").Value = -54845.8

$wsARM.Range("H139").Value = 147999.6
$wsARM.Range("J139").Value = 113333
$wsARM.Range("L139").Value = 113333
$wsARM.Range("N139").Value = -123613

$wsBSM.Range("H105").Value = 2129.4211
$wsBSM.Range("I105").Value = 1821
$wsBSM.Range("K105").Value = 1821
$wsBSM.Range("M105").Value = -74

$wsBSM.Range("H107").Value = 2288.182
$wsBSM.Range("I107").Value = 2187.652
$wsBSM.Range("J107").Value = 2519.4
$wsBSM.Range("K107").Value = 2187.652
$wsBSM.Range("L107").Value = 2519.4
$wsBSM.Range("M107").Value = -267.652
$wsBSM.Range("N107").Value = -6359.4

$wsCRP.Range("H7").Value = 144.66667
$wsCRP.Range("I7").Value = 72.40000000000001
$wsCRP.Range("J7").Value = 506
$wsCRP.Range("K7").Value = 72.40000000000001
$wsCRP.Range("L7").Value = 506
$wsCRP.Range("M7").Value = 40.59999999999999
$wsCRP.Range("N7").Value = -732

$wsCRP.Range("H10").Value = 1173.1875
$wsCRP.Range("I10").Value = 1078.9166
$wsCRP.Range("J10").Value = 1456
$wsCRP.Range("K10").Value = 1078.9166
$wsCRP.Range("L10").Value = 1456
$wsCRP.Range("M10").Value = -939.9166
$wsCRP.Range("N10").Value = -1734

$wsCRP.Range("H31").Value = 2941.7437
$wsCRP.Range("I31").Value = 2311.087
$wsCRP.Range("J31").Value = 3848.3125
$wsCRP.Range("K31").Value = 2311.087
$wsCRP.Range("L31").Value = 3848.3125
$wsCRP.Range("M31").Value = -2016.087
$wsCRP.Range("N31").Value = -4438.3125

$wsCRP.Range("H34").Value = 2941.7437
$wsCRP.Range("I34").Value = 2311.087
$wsCRP.Range("J34").Value = 3848.3125
$wsCRP.Range("K34").Value = 2311.087
$wsCRP.Range("L34").Value = 3848.3125
$wsCRP.Range("M34").Value = -2109.087
$wsCRP.Range("N34").Value = -4252.3125

$wsCRP.Range("H105").Value = 2511.6875
$wsCRP.Range("I105").Value = 2245.2307
$wsCRP.Range("K105").Value = 2245.2307
$wsCRP.Range("M105").Value = -498.2307000000001

$wsCRP.Range("H122").Value = 1329.6957
$wsCRP.Range("I122").Value = 1378.0526
$wsCRP.Range("K122").Value = 4134.1578
$wsCRP.Range("M122").Value = -1684.1578

$wsCRP.Range("H134").Value = 2892
$wsCRP.Range("I134").Value = 2304.2727
$wsCRP.Range("K134").Value = 6912.8181
$wsCRP.Range("M134").Value = -4377.8181

$wsCUL.Range("H44").Value = 2918.4443
$wsCUL.Range("I44").Value = 909.4286
$wsCUL.Range("K44").Value = 2728.2858
$wsCUL.Range("M44").Value = -2330.2858

$wsCUL.Range("H97").Value = 944.63635
$wsCUL.Range("I97").Value = 627.8570999999999
$wsCUL.Range("J97").Value = 1499
$wsCUL.Range("K97").Value = 1883.5713
$wsCUL.Range("L97").Value = 4497
$wsCUL.Range("M97").Value = -1387.5713
$wsCUL.Range("N97").Value = -5489

$wsGSM.Range("H11").Value = 14166667
$wsGSM.Range("I11").Value = 16250000
$wsGSM.Range("K11").Value = 16250000
$wsGSM.Range("M11").Value = -16249861

$wsGSM.Range("H70").Value = 4500.1626
$wsGSM.Range("I70").Value = 4502
$wsGSM.Range("K70").Value = 4502
$wsGSM.Range("M70").Value = -4232

$wsGSM.Range("H73").Value = 4500.1626
$wsGSM.Range("I73").Value = 4502
$wsGSM.Range("K73").Value = 4502
$wsGSM.Range("M73").Value = -3566

$wsGSM.Range("H126").Value = 4481.3335
$wsGSM.Range("I126").Value = 3678.8333
$wsGSM.Range("K126").Value = 11036.4999
$wsGSM.Range("M126").Value = -8566.499899999999

$wsGSM.Range("H132").Value = 3892.889
$wsGSM.Range("I132").Value = 3858.2856
$wsGSM.Range("K132").Value = 11574.8568
$wsGSM.Range("M132").Value = -9044.856800000001

$wsLTW.Range("H7").Value = 6507.9653
$wsLTW.Range("I7").Value = 7337.8696
$wsLTW.Range("J7").Value = 3326.6667
$wsLTW.Range("K7").Value = 7337.8696
$wsLTW.Range("L7").Value = 3326.6667
$wsLTW.Range("M7").Value = -7225.8696
$wsLTW.Range("N7").Value = -3550.6667

$wsLTW.Range("H40").Value = 4973.5
$wsLTW.Range("I40").Value = 4692.727
$wsLTW.Range("J40").Value = 6003
$wsLTW.Range("K40").Value = 4692.727
$wsLTW.Range("L40").Value = 6003
$wsLTW.Range("M40").Value = -4556.727
$wsLTW.Range("N40").Value = -6275

$wsLTW.Range("H46").Value = 2785
$wsLTW.Range("I46").Value = 750
$wsLTW.Range("K46").Value = 750
$wsLTW.Range("M46").Value = -562

$wsLTW.Range("H126").Value = 6507.9653
$wsLTW.Range("I126").Value = 7337.8696
$wsLTW.Range("J126").Value = 3326.6667
$wsLTW.Range("K126").Value = 22013.6088
$wsLTW.Range("L126").Value = 9980.000100000001
$wsLTW.Range("M126").Value = -19543.6088
$wsLTW.Range("N126").Value = -14920.0001

$wsWVR.Range("H81").Value = 2488.3333
$wsWVR.Range("I81").Value = 1934.7142
$wsWVR.Range("J81").Value = 2840.6365
$wsWVR.Range("K81").Value = 3869.4284
$wsWVR.Range("L81").Value = 5681.273
$wsWVR.Range("M81").Value = -2808.4284
$wsWVR.Range("N81").Value = -7803.273

$wsWVR.Range("H84").Value = 2488.3333
$wsWVR.Range("I84").Value = 1934.7142
$wsWVR.Range("J84").Value = 2840.6365
$wsWVR.Range("K84").Value = 19347.142
$wsWVR.Range("L84").Value = 28406.365
$wsWVR.Range("M84").Value = -14043.142
$wsWVR.Range("N84").Value = -39014.36500000001

$wsWVR.Range("H122").Value = 3440.1853
$wsWVR.Range("I122").Value = 1575.2106
$wsWVR.Range("J122").Value = 7869.5
$wsWVR.Range("K122").Value = 4725.6318
$wsWVR.Range("L122").Value = 23608.5
$wsWVR.Range("M122").Value = -2275.6318
$wsWVR.Range("N122").Value = -28508.5

$wsWVR.Range("H132").Value = 2341.164
$wsWVR.Range("I132").Value = 2304.1897
$wsWVR.Range("J132").Value = 2579.4443
$wsWVR.Range("K132").Value = 6912.5691
$wsWVR.Range("L132").Value = 7738.3329
$wsWVR.Range("M132").Value = -4382.5691
$wsWVR.Range("N132").Value = -12798.3329
